$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each self-assessment indicator group has a row-range in column G (already present,
# merged across the group). We add matching merged cells in columns H (Completamente
# Logrado), I (Logrado), J (Logro incipiente) and K (No logrado), mark the achieved
# level with an "x", and draw a thin box border around every one of the 4 new cells
# (Excel naturally fragments that border across the merged rows).

$groups = @(
    @{ Rows = @(6);        Col = "H" },
    @{ Rows = @(7,8,9);    Col = "H" },
    @{ Rows = @(10,11,12); Col = "H" },
    @{ Rows = @(13,14,15); Col = "H" },
    @{ Rows = @(16);       Col = "H" },
    @{ Rows = @(17);       Col = "H" },
    @{ Rows = @(18,19);    Col = "H" },
    @{ Rows = @(20);       Col = "I" },
    @{ Rows = @(21,22,23); Col = "I" },
    @{ Rows = @(24);       Col = "H" },
    @{ Rows = @(25);       Col = "H" },
    @{ Rows = @(26,27,28); Col = "K" }
)

$columns = @("H", "I", "J", "K")

foreach ($g in $groups) {
    $rows = $g.Rows
    $firstRow = $rows[0]
    $lastRow = $rows[$rows.Count - 1]
    $xCol = $g.Col

    foreach ($col in $columns) {
        $addr = "$col$firstRow`:$col$lastRow"
        $rng = $ws.Range($addr)
        if ($rows.Count -gt 1) {
            $rng.Merge() | Out-Null
        }
        $rng.HorizontalAlignment = -4108
        $rng.VerticalAlignment = -4108
        $rng.BorderAround(1, 2) | Out-Null
        if ($col -eq $xCol) {
            $rng.Value2 = "x"
        }
    }
}

# Total row (row 29): boxed but empty cells, no merge across rows.
foreach ($col in $columns) {
    $rng = $ws.Range("$col`29")
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4108
    $rng.BorderAround(1, 2) | Out-Null
}

$ws.Range("B25").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 25
$ws.Range("K25").Select() | Out-Null
